$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = ' Churchill-laan 290 2, 1078 GB Amsterdam Verkocht onder voorbehoud '
$ws.Range("C2").Value = 870000
$ws.Range("D2").Value = 98
$ws.Range("E2").Value = 'C'
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 1926
$ws.Range("J2").Value = 'Goed'
$ws.Range("K2").Value = 'Goed'
$ws.Range("L2").Value = 0.7941636363636363

# Row 3
$ws.Range("B3").Value = ' Churchill-laan 157 2, 1078 DV Amsterdam Verkocht Width'
$ws.Range("C3").Value = 675000
$ws.Range("D3").Value = 94
$ws.Range("E3").Value = 'C'
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1928
$ws.Range("J3").Value = 'Goed'
$ws.Range("K3").Value = 'Goed'
$ws.Range("L3").Value = 0.7797636363636363

# Row 4
$ws.Range("B4").Value = ' Churchill-laan 153 4, 1078 DT Amsterdam Verkocht Width'
$ws.Range("C4").Value = 795000
$ws.Range("D4").Value = 89
$ws.Range("E4").Value = 'C'
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1928
$ws.Range("J4").Value = 'Goed'
$ws.Range("K4").Value = 'Goed'
$ws.Range("L4").Value = 0.7617636363636363

# Row 5
$ws.Range("B5").Value = ' Scheldestraat 29 1, 1078 GE Amsterdam Verkocht Width'
$ws.Range("C5").Value = 850000
$ws.Range("D5").Value = 95
$ws.Range("E5").Value = 'C'
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1928
$ws.Range("J5").Value = 'Goed'
$ws.Range("K5").Value = 'Goed'
$ws.Range("L5").Value = 0.7547922077922078

# Row 6
$ws.Range("B6").Value = ' Churchill-laan 119 4, 1078 DN Amsterdam Verkocht Width'
$ws.Range("C6").Value = 675000
$ws.Range("D6").Value = 87
$ws.Range("E6").Value = 'C'
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1928
$ws.Range("J6").Value = 'Goed'
$ws.Range("K6").Value = 'Goed'
$ws.Range("L6").Value = 0.7545636363636363

# Row 7
$ws.Range("B7").Value = ' Churchill-laan 282 3, 1078 GB Amsterdam Verkocht Width'
$ws.Range("C7").Value = 795000
$ws.Range("D7").Value = 97
$ws.Range("E7").Value = 'D'
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1929
$ws.Range("J7").Value = 'Goed'
$ws.Range("K7").Value = 'Goed'
$ws.Range("L7").Value = 0.7537454545454545

# Row 8
$ws.Range("B8").Value = ' Scheldestraat 29 4, 1078 GE Amsterdam Verkocht Width'
$ws.Range("C8").Value = 700000
$ws.Range("D8").Value = 90
$ws.Range("E8").Value = 'B'
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1928
$ws.Range("J8").Value = 'Goed'
$ws.Range("K8").Value = 'Goed'
$ws.Range("L8").Value = 0.7486103896103896

# Row 9
$ws.Range("B9").Value = ' Rooseveltlaan 168 2, 1078 NT Amsterdam Verkocht Width'
$ws.Range("C9").Value = 875000
$ws.Range("D9").Value = 101
$ws.Range("E9").Value = 'C'
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1934
$ws.Range("J9").Value = 'Goed'
$ws.Range("K9").Value = 'Goed'
$ws.Range("L9").Value = 0.7477636363636363

# Row 10
$ws.Range("B10").Value = ' Rooseveltlaan 224 III, 1078 NX Amsterdam Verkocht Width'
$ws.Range("C10").Value = 800000
$ws.Range("D10").Value = 98
$ws.Range("E10").Value = 'C'
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1934
$ws.Range("J10").Value = 'Uitstekend'
$ws.Range("K10").Value = 'Uitstekend'
$ws.Range("L10").Value = 0.7441636363636364

# Row 11
$ws.Range("B11").Value = ' Rooseveltlaan 230 3, 1078 NX Amsterdam Verkocht Width'
$ws.Range("C11").Value = 895000
$ws.Range("D11").Value = 98
$ws.Range("E11").Value = 'A'
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1934
$ws.Range("J11").Value = 'Goed'
$ws.Range("K11").Value = 'Goed'
$ws.Range("L11").Value = 0.7428000000000001

# Row 12
$ws.Range("B12").Value = ' Churchill-laan 248 1, 1078 EZ Amsterdam Verkocht Width'
$ws.Range("C12").Value = 750000
$ws.Range("D12").Value = 105
$ws.Range("E12").Value = 'Unknown'
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 1928
$ws.Range("J12").Value = 'Matig'
$ws.Range("K12").Value = 'Matig'
$ws.Range("L12").Value = 0.742

# Row 13
$ws.Range("B13").Value = ' Churchill-laan 161 H, 1078 DV Amsterdam Verkocht Width'
$ws.Range("C13").Value = 995000
$ws.Range("D13").Value = 110
$ws.Range("E13").Value = 'D'
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1928
$ws.Range("I13").Value = $true
$ws.Range("J13").Value = 'Goed'
$ws.Range("K13").Value = 'Goed'
$ws.Range("L13").Value = 0.7385454545454547

# Row 14
$ws.Range("B14").Value = ' Churchill-laan 59 B, 1078 DH Amsterdam Verkocht Width'
$ws.Range("C14").Value = 715000
$ws.Range("D14").Value = 87
$ws.Range("E14").Value = 'C'
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1927
$ws.Range("J14").Value = 'Uitstekend'
$ws.Range("K14").Value = 'Uitstekend'
$ws.Range("L14").Value = 0.7295636363636364

# Row 15
$ws.Range("B15").Value = ' Rooseveltlaan 122 1, 1078 NP Amsterdam Verkocht Width'
$ws.Range("C15").Value = 899000
$ws.Range("D15").Value = 101
$ws.Range("E15").Value = 'C'
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 1930
$ws.Range("J15").Value = 'Uitstekend'
$ws.Range("K15").Value = 'Uitstekend'
$ws.Range("L15").Value = 0.7227636363636363

# Row 16
$ws.Range("B16").Value = ' Amstelkade 168 2, 1078 AZ Amsterdam Verkocht Width'
$ws.Range("C16").Value = 850000
$ws.Range("D16").Value = 102
$ws.Range("E16").Value = 'D'
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1936
$ws.Range("J16").Value = 'Goed'
$ws.Range("K16").Value = 'Goed'
$ws.Range("L16").Value = 0.7180597402597403
